# TPS Report weekly update (Sept 21 - Sept 27)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 12: Installation Documentation ---
# Risk went down (10 -> 4), Time Actual went up (4 Hours -> 7 Hours), % Complete up (0.5 -> 0.6)
$ws.Range("C12").Value = 4
$ws.Range("D12").Value = "7 Hours"
$ws.Range("E12").Value = 0.6

# --- Row 13: Add Humidity modification ---
# Risk went down (10 -> 6), Time Actual went up (1.5 Hours -> 3 Hours)
$ws.Range("C13").Value = 6
$ws.Range("D13").Value = "3 Hours"

# --- Row 14 (new): Create Ruby Version of Site ---
$ws.Range("A14").Value = "Create Ruby Version of Site"
$ws.Range("B14").Value = "45 Hours"
$ws.Range("C14").Value = 10
$ws.Range("D14").Value = "12 Hours"
$ws.Range("E14").Value = 0.15
$ws.Range("E14").NumberFormat = "0.00%"
$ws.Range("F14").Value = "All"

# Move / update the active selection to reflect where the user ended up editing
$null = $ws.Range("F20").Select()
